$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 861
$ws.Cells.Item(5, 6).Value = 1177
$ws.Cells.Item(6, 6).Value = 57
$ws.Cells.Item(7, 6).Value = 4267
$ws.Cells.Item(8, 6).Value = 2571
$ws.Cells.Item(10, 6).Value = 2457
$ws.Cells.Item(12, 6).Value = 1946
$ws.Cells.Item(13, 6).Value = 489
$ws.Cells.Item(14, 6).Value = 1646
$ws.Cells.Item(15, 6).Value = 653
$ws.Cells.Item(16, 6).Value = 15
$ws.Cells.Item(18, 6).Value = 311
$ws.Cells.Item(20, 6).Value = 268
$ws.Cells.Item(21, 6).Value = 70
$ws.Cells.Item(22, 6).Value = 11
$ws.Cells.Item(23, 6).Value = 458
$ws.Cells.Item(26, 6).Value = 510
$ws.Cells.Item(27, 6).Value = 682
$ws.Cells.Item(28, 6).Value = 92
$ws.Cells.Item(30, 6).Value = 384
$ws.Cells.Item(31, 6).Value = 41
$ws.Cells.Item(32, 6).Value = 1612
$ws.Cells.Item(33, 6).Value = 942
$ws.Cells.Item(34, 6).Value = 72
$ws.Cells.Item(36, 6).Value = 1026
$ws.Cells.Item(37, 6).Value = 1994
$ws.Cells.Item(38, 6).Value = 244
$ws.Cells.Item(40, 6).Value = 527
$ws.Cells.Item(41, 6).Value = 82
$ws.Cells.Item(42, 6).Value = 15
$ws.Cells.Item(43, 6).Value = 615
$ws.Cells.Item(44, 6).Value = 1279
$ws.Cells.Item(45, 6).Value = 69
$ws.Cells.Item(47, 6).Value = 418
$ws.Cells.Item(48, 6).Value = 60

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 8
$ws.Cells.Item(9, 6).Value = 10

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 861
$ws.Cells.Item(3, 6).Value = 1177
$ws.Cells.Item(4, 6).Value = 8
$ws.Cells.Item(5, 6).Value = 57
$ws.Cells.Item(6, 6).Value = 4267
$ws.Cells.Item(7, 6).Value = 2571
$ws.Cells.Item(8, 6).Value = 2457
$ws.Cells.Item(9, 6).Value = 1946
$ws.Cells.Item(10, 6).Value = 1646
$ws.Cells.Item(12, 6).Value = 653
$ws.Cells.Item(13, 6).Value = 15
$ws.Cells.Item(15, 6).Value = 311
$ws.Cells.Item(17, 6).Value = 268
$ws.Cells.Item(18, 6).Value = 70
$ws.Cells.Item(19, 6).Value = 458
$ws.Cells.Item(22, 6).Value = 510
$ws.Cells.Item(23, 6).Value = 682
$ws.Cells.Item(24, 6).Value = 92
$ws.Cells.Item(29, 6).Value = 384
$ws.Cells.Item(30, 6).Value = 1612
$ws.Cells.Item(31, 6).Value = 942
$ws.Cells.Item(32, 6).Value = 72
$ws.Cells.Item(35, 6).Value = 1026
$ws.Cells.Item(36, 6).Value = 1994
$ws.Cells.Item(37, 6).Value = 244
$ws.Cells.Item(39, 6).Value = 10
$ws.Cells.Item(41, 6).Value = 527
$ws.Cells.Item(42, 6).Value = 82
$ws.Cells.Item(43, 6).Value = 15
$ws.Cells.Item(44, 6).Value = 615
$ws.Cells.Item(45, 6).Value = 1279
$ws.Cells.Item(46, 6).Value = 69
$ws.Cells.Item(47, 6).Value = 418
$ws.Cells.Item(48, 6).Value = 60
